$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("StaffSalariesDB")

# Insert a new row at row 87 (pushing existing rows 87..180 down to 88..181)
$ws.Rows.Item(87).Insert()

$ws.Range("A87").Value = 296
$ws.Range("B87").Value = "Kiribati"
$ws.Range("C87").Value = "KIR"
$ws.Range("D87").Value = 4238.7462939999996
$ws.Range("E87").Value = 5510.370183
$ws.Range("F87").Value = 7417.8060150000001
$ws.Range("G87").Value = 9007.3358759999992
$ws.Range("H87").Value = 13511.00381

# Match the author's final view state: scrolled down with J82 selected
$excel.ActiveWindow.ScrollRow = 58
[void]$ws.Range("J82").Select()
